$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top to hold the header labels,
# shifting the existing data rows down.
$ws.Rows.Item(1).Insert()

# Populate the new header row with column labels.
$ws.Range("A1").Value = "matric"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "department"
$ws.Range("D1").Value = "level"
$ws.Range("E1").Value = "ca"

# Update selection to match target state.
$ws.Range("E1").Select()
